$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data table (header unchanged).
# Two new parameter rows (beta1, beta2) inserted at the top of the
# parameter list, shifting the remaining parameters (gamma1, gamma2, nu,
# rho, phi1, phi2) down by two rows, with new fitted values.

$data = @(
    @(17, "beta1",  0.97, 0.97, "( 0.84 , 1.06 )", 127.9,  1.03),
    @(17, "beta2",  0.58, 0.51, "( 0.04 , 1.43 )", 797.82, 1),
    @(17, "gamma1", 0.25, 0.25, "( 0.15 , 0.34 )", 203.88, 1.02),
    @(17, "gamma2", 0.25, 0.25, "( 0.15 , 0.35 )", 450.36, 1),
    @(17, "nu",     0,    0,    "( 0 , 0 )",       342.94, 1.01),
    @(17, "rho",    0.34, 0.33, "( 0.06 , 0.71 )", 584.59, 1),
    @(17, "phi1",   0.49, 0.49, "( 0.34 , 0.7 )",  395.88, 1.01),
    @(17, "phi2",   0,    0,    "( 0 , 0 )",       631.57, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
